# Apply feedback/validation review edits described in commit:
# "Validaciones DG, CL, PL - Se están corroborando los ajustes de
#  validaciones por parte de la retroalimentación recibida"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------
# 1) Mark the "Estatus" (status) column E -- and, for a few rows, the
#    mirrored column H -- as "Completa" for the rows that were reviewed.
# ---------------------------------------------------------------------
$completaRowsE = @(6,7,8,9,10,11,12,13,14,15,16,20,21,22,23,24,25,26,27,28,29,30,31,32)
foreach ($r in $completaRowsE) {
    $ws.Range("E$r").Value = "Completa"
}

$completaRowsH = @(6,11,12)
foreach ($r in $completaRowsH) {
    $ws.Range("H$r").Value = "Completa"
}

# ---------------------------------------------------------------------
# 2) For the rows whose review turned out incomplete, note it in the
#    adjoining "Comentarios" column F.
# ---------------------------------------------------------------------
$faltaRowsF = @(23,24,25,26,31,32)
foreach ($r in $faltaRowsF) {
    $ws.Range("F$r").Value = "falta una parte"
}

# ---------------------------------------------------------------------
# 3) Re-apply / extend the conditional formatting that colors column E
#    according to its status, stretching it to cover the whole column
#    (E3:E517) instead of just the originally filled rows (E3:E17).
# ---------------------------------------------------------------------
$cfRange = $ws.Range("E3:E517")
$cfRange.FormatConditions.Delete()

$fcPendiente = $cfRange.FormatConditions.Add(1, 3, '"Pendiente"')
$fcPendiente.Font.Color = 22428
$fcPendiente.Interior.Color = 10284031

$fcCompleta = $cfRange.FormatConditions.Add(1, 3, '"Completa"')
$fcCompleta.Font.Color = 24832
$fcCompleta.Interior.Color = 13561798

$fcCompleta.Priority = 1
$fcPendiente.Priority = 2

# ---------------------------------------------------------------------
# 4) Update the sheet view: scroll the frozen pane down and select the
#    cell the reviewer ended up on.
# ---------------------------------------------------------------------
$ws.Range("F34").Select()

$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
$win.Width = 20730
$win.Height = 11310
